$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column (D2:D51) as Text first so that numeric-looking
# strings (e.g. "1.001", "0.1240") are stored verbatim instead of being
# auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.109.18"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.667.16"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "209.27"
$ws.Range("E5").Value = "  -3.84%  "
$ws.Range("D6").Value = "0.5247"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "0.2614"
$ws.Range("E8").Value = "  -3.86%  "
$ws.Range("D9").Value = "0.06292"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "21.10"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("D11").Value = "0.07529"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "1.667.38"
$ws.Range("E12").Value = "  -6.95%  "
$ws.Range("D13").Value = "4.429"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "0.5504"
$ws.Range("E14").Value = "  -4.92%  "
$ws.Range("D15").Value = "66.39"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "0.000007921"
$ws.Range("E16").Value = "  -5.21%  "
$ws.Range("D17").Value = "26.129.66"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "4.709"
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("D20").Value = "186.25"
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").Value = "10.25"
$ws.Range("E21").Value = "  -5.45%  "
$ws.Range("D22").Value = "6.164"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "149.48"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "0.1240"
$ws.Range("E25").Value = "  -3.69%  "
$ws.Range("D26").Value = "7.448"
$ws.Range("E26").Value = "  -5.08%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "0.06341"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("D29").Value = "1.353"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").Value = "1.274"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").Value = "3.489"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").Value = "3.409"
$ws.Range("E32").Value = "  -4.81%  "
$ws.Range("D33").Value = "1.633"
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("D34").Value = "1.001"
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").Value = "2.407"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "0.6010"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").Value = "2.728"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "1.105.65"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "6.090"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "0.01613"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "0.8688"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "99.83"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "1.817.89"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "0.00000000107"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").Value = "55.24"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "8.025"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "0.05229"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "0.4245"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "5.926"
$ws.Range("E51").Value = "  -2.04%  "

# Revert the Price column formatting/style back to the default so the
# workbook styles are unchanged from before the edit.
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"

